# Generate Report for Handback
# Update the generated timestamps recorded on the handback status report.

$wb = $excel.ActiveWorkbook

# "Overview" sheet - Latest HO Xliff Generate Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-23 21:08:43"

# "zh-cn" sheet - Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-23 21:08:39"
$wsZhCn.Range("K2").Value = "2016-08-23 21:08:56"

# "de-de" sheet - Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-23 21:08:43"
$wsDeDe.Range("K2").Value = "2016-08-23 21:09:10"
